$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 66.47695399999999
$ws.Cells.Item(2, 8).Value = 199.430862
$ws.Cells.Item(2, 9).Value = 0.04311983106164722
$ws.Cells.Item(2, 10).Value = 0.04311983106164721
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4910443333333334
$ws.Cells.Item(2, 14).Value = 1.473133
$ws.Cells.Item(2, 15).Value = 0.7844104380534107
$ws.Cells.Item(2, 16).Value = 0.7844104380534107
$ws.Cells.Item(2, 17).Value = 32.64313155896066
$ws.Cells.Item(2, 18).Value = 293.788184030646
$ws.Cells.Item(2, 19).Value = 0.03382364557185576
$ws.Cells.Item(2, 20).Value = 0.03382364557185576

$ws.Cells.Item(3, 7).Value = 66.47695399999999
$ws.Cells.Item(3, 8).Value = 199.430862
$ws.Cells.Item(3, 9).Value = 0.04311983106164722
$ws.Cells.Item(3, 10).Value = 0.04311983106164721
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.13496
$ws.Cells.Item(3, 14).Value = 0.40488
$ws.Cells.Item(3, 15).Value = 0.2155895619465893
$ws.Cells.Item(3, 16).Value = 0.2155895619465893
$ws.Cells.Item(3, 17).Value = 8.971729711839998
$ws.Cells.Item(3, 18).Value = 80.74556740656
$ws.Cells.Item(3, 19).Value = 0.009296185489791458
$ws.Cells.Item(3, 20).Value = 0.009296185489791458

$ws.Cells.Item(4, 9).Value = 0.8830494168872806
$ws.Cells.Item(4, 10).Value = 0.8830494168872804
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4910443333333334
$ws.Cells.Item(4, 14).Value = 1.473133
$ws.Cells.Item(4, 15).Value = 0.7844104380534107
$ws.Cells.Item(4, 16).Value = 0.7844104380534107
$ws.Cells.Item(4, 17).Value = 668.4974773510589
$ws.Cells.Item(4, 18).Value = 6016.477296159531
$ws.Cells.Item(4, 19).Value = 0.6926731799233606
$ws.Cells.Item(4, 20).Value = 0.6926731799233605

$ws.Cells.Item(5, 9).Value = 0.8830494168872806
$ws.Cells.Item(5, 10).Value = 0.8830494168872804
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.13496
$ws.Cells.Item(5, 14).Value = 0.40488
$ws.Cells.Item(5, 15).Value = 0.2155895619465893
$ws.Cells.Item(5, 16).Value = 0.2155895619465893
$ws.Cells.Item(5, 17).Value = 183.73171915224
$ws.Cells.Item(5, 18).Value = 1653.58547237016
$ws.Cells.Item(5, 19).Value = 0.1903762369639199
$ws.Cells.Item(5, 20).Value = 0.1903762369639199

$ws.Cells.Item(6, 7).Value = 44.831112
$ws.Cells.Item(6, 8).Value = 134.493336
$ws.Cells.Item(6, 9).Value = 0.02907940059566787
$ws.Cells.Item(6, 10).Value = 0.02907940059566786
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.4910443333333334
$ws.Cells.Item(6, 14).Value = 1.473133
$ws.Cells.Item(6, 15).Value = 0.7844104380534107
$ws.Cells.Item(6, 16).Value = 0.7844104380534107
$ws.Cells.Item(6, 17).Value = 22.014063504632
$ws.Cells.Item(6, 18).Value = 198.126571541688
$ws.Cells.Item(6, 19).Value = 0.02281018535957845
$ws.Cells.Item(6, 20).Value = 0.02281018535957844

$ws.Cells.Item(7, 7).Value = 44.831112
$ws.Cells.Item(7, 8).Value = 134.493336
$ws.Cells.Item(7, 9).Value = 0.02907940059566787
$ws.Cells.Item(7, 10).Value = 0.02907940059566786
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.13496
$ws.Cells.Item(7, 14).Value = 0.40488
$ws.Cells.Item(7, 15).Value = 0.2155895619465893
$ws.Cells.Item(7, 16).Value = 0.2155895619465893
$ws.Cells.Item(7, 17).Value = 6.050406875519999
$ws.Cells.Item(7, 18).Value = 54.45366187968001
$ws.Cells.Item(7, 19).Value = 0.006269215236089423
$ws.Cells.Item(7, 20).Value = 0.006269215236089422

$ws.Cells.Item(8, 7).Value = 52.83062100000001
$ws.Cells.Item(8, 8).Value = 158.491863
$ws.Cells.Item(8, 9).Value = 0.0342682285413064
$ws.Cells.Item(8, 10).Value = 0.03426822854130639
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.4910443333333334
$ws.Cells.Item(8, 14).Value = 1.473133
$ws.Cells.Item(8, 15).Value = 0.7844104380534107
$ws.Cells.Item(8, 16).Value = 0.7844104380534107
$ws.Cells.Item(8, 17).Value = 25.94217706853101
$ws.Cells.Item(8, 18).Value = 233.479593616779
$ws.Cells.Item(8, 19).Value = 0.02688035616140054
$ws.Cells.Item(8, 20).Value = 0.02688035616140054

$ws.Cells.Item(9, 7).Value = 52.83062100000001
$ws.Cells.Item(9, 8).Value = 158.491863
$ws.Cells.Item(9, 9).Value = 0.0342682285413064
$ws.Cells.Item(9, 10).Value = 0.03426822854130639
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.13496
$ws.Cells.Item(9, 14).Value = 0.40488
$ws.Cells.Item(9, 15).Value = 0.2155895619465893
$ws.Cells.Item(9, 16).Value = 0.2155895619465893
$ws.Cells.Item(9, 17).Value = 7.130020610160001
$ws.Cells.Item(9, 18).Value = 64.17018549144001
$ws.Cells.Item(9, 19).Value = 0.007387872379905854
$ws.Cells.Item(9, 20).Value = 0.007387872379905853

$ws.Cells.Item(10, 7).Value = 16.16161433333333
$ws.Cells.Item(10, 8).Value = 48.484843
$ws.Cells.Item(10, 9).Value = 0.01048312291409786
$ws.Cells.Item(10, 10).Value = 0.01048312291409786
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4910443333333334
$ws.Cells.Item(10, 14).Value = 1.473133
$ws.Cells.Item(10, 15).Value = 0.7844104380534107
$ws.Cells.Item(10, 16).Value = 0.7844104380534107
$ws.Cells.Item(10, 17).Value = 7.936069135902112
$ws.Cells.Item(10, 18).Value = 71.424622223119
$ws.Cells.Item(10, 19).Value = 0.008223071037215254
$ws.Cells.Item(10, 20).Value = 0.008223071037215252

$ws.Cells.Item(11, 7).Value = 16.16161433333333
$ws.Cells.Item(11, 8).Value = 48.484843
$ws.Cells.Item(11, 9).Value = 0.01048312291409786
$ws.Cells.Item(11, 10).Value = 0.01048312291409786
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.13496
$ws.Cells.Item(11, 14).Value = 0.40488
$ws.Cells.Item(11, 15).Value = 0.2155895619465893
$ws.Cells.Item(11, 16).Value = 0.2155895619465893
$ws.Cells.Item(11, 17).Value = 2.181171470426666
$ws.Cells.Item(11, 18).Value = 19.63054323384
$ws.Cells.Item(11, 19).Value = 0.002260051876882611
$ws.Cells.Item(11, 20).Value = 0.002260051876882611
